# Experiment 01 - LLMs in TCM to MBT initial commit
# Applies the shared-string content corrections to the "Test Suite" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# --- Header: System name correction ---
$ws.Range("B1").Value = "GTI-Competências"

# --- TC2 block (rows 16-23): step 2 expected result ---
$ws.Range("D21").Value = "SYSTEM alerta que o nome de usuário e/ou senha estão incorretos"

# --- TC3 block (rows 26-33): test id + step 2 data/expected result ---
$ws.Range("B26").Value = "TC3"
$ws.Range("B31").Value = "Usuário do Sistema preenche os campos e clica no botão entrar"
$ws.Range("D31").Value = "SYSTEM alerta que o CAS (sistema de autorização login-senha) está fora do ar"

# --- TC4 block (rows 36-43): test id + step 2 expected result ---
$ws.Range("B36").Value = "TC4"
$ws.Range("D41").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissões de acesso e escrita) está fora do ar"

# --- TC5 block (rows 46-53): step 2/3 data + step 2 expected result ---
$ws.Range("B51").Value = "Usuário do Sistema seleciona um nome de usuário sugerido, digita a senha e clica no botão entrar"
$ws.Range("D51").Value = "SYSTEM alerta que o CAS (sistema de autorização login-senha) está fora do ar"
$ws.Range("B52").Value = "Usuário do Sistema seleciona um nome de usuário sugerido, digita a senha e clica no botão entrar"

# --- TC6 block (rows 56-63): step 2 expected result + step 3 data ---
$ws.Range("D61").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissões de acesso e escrita) está fora do ar"
$ws.Range("B62").Value = "Usuário do Sistema preenche os campos e clica no botão entrar"
